$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

$ws.Range("A2").Value = "a"
$ws.Range("B2").Value = "incomplete"
$ws.Range("C2").Value = "2025-01-10 15:57:00.352011"
$ws.Range("D2").Value = "N/A"

$ws.Range("A3").Value = "a"
$ws.Range("B3").Value = "incomplete"
$ws.Range("C3").Value = "2025-01-10 15:57:00.499081"
$ws.Range("D3").Value = "N/A"
